$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2, shifting the existing data down.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with "Not Applicable" across A2:D2.
$ws.Range("A2:D2").Value = "Not Applicable"

# The inserted row inherits the bold header formatting; the new row should
# remain regular (unstyled) like the other data rows.
$ws.Range("A2:D2").Font.Bold = $false

# Update the selection to match the new state (A2:D2 selected).
$ws.Range("A2:D2").Select()

# Adjust column widths to fit the new, wider content.
$ws.Columns.Item(1).ColumnWidth = 13.25
$ws.Columns.Item(2).ColumnWidth = 12.1
$ws.Columns.Item(3).ColumnWidth = 13.25
$ws.Columns.Item(4).ColumnWidth = 12.1
